$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 3 (full format + values) into rows 4, 5 and 6.
foreach ($row in 4,5,6) {
    $ws.Range("A3:T3").Copy()
    $ws.Range("A${row}:T${row}").PasteSpecial(-4122)
    $ws.Range("A3:T3").Copy()
    $ws.Range("A${row}:T${row}").PasteSpecial(-4163)
    $ws.Rows.Item($row).RowHeight = 15
}

# Names first (column C), matching the order the new strings were authored in.
$ws.Range("C4").Value = "Peeyush Vatsi"
$ws.Range("C5").Value = "Harsh"
$ws.Range("C6").Value = "Bagdadi"

# Then parent names (column I).
$ws.Range("I4").Value = "abc"
$ws.Range("I5").Value = "def"
$ws.Range("I6").Value = "ghi"

# Student IDs.
$ws.Range("B4").Value = 19105008
$ws.Range("B5").Value = 19105030
$ws.Range("B6").Value = 19105015

# Leave flags flip relative to row 3 (applied=true, approved=false).
foreach ($row in 4,5,6) {
    $ws.Range("N$row").Value = $true
    $ws.Range("O$row").Value = $false
}

$ws.Range("O6").Select()
